# Updates weekly Fruta/Hortaliza price data for
# "Fruta, Terminal La Palmera de La Serena - Coco" (commit: "Fruta / hortaliza, semanal").
# Each data row (2..41) gets its Fecha (D), Volumen (M), Precio minimo (N),
# Precio maximo (O), Precio promedio ponderado (P) and Precio $/Kg (S)
# columns refreshed to the latest weekly snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D=44462; M=100; N=19500; O=20000; P=19750; S=988},
    @{Row=3; D=44333; M=100; N=19500; O=20000; P=19750; S=988},
    @{Row=4; D=44364; M=140; N=20000; O=21000; P=20500; S=1025},
    @{Row=5; D=44445; M=160; N=20000; O=21000; P=20500; S=1025},
    @{Row=6; D=44434; M=100; N=20000; O=21000; P=20500; S=1025},
    @{Row=7; D=44350; M=160; N=19000; O=20000; P=19500; S=975},
    @{Row=8; D=44442; M=140; N=20000; O=21000; P=20500; S=1025},
    @{Row=9; D=44467; M=200; N=20000; O=21000; P=20500; S=1025},
    @{Row=10; D=44428; M=100; N=20000; O=21000; P=20500; S=1025},
    @{Row=11; D=44474; M=200; N=19000; O=20000; P=19500; S=975},
    @{Row=12; D=44418; M=200; N=20000; O=21000; P=20500; S=1025},
    @{Row=13; D=44879; M=100; N=28000; O=30000; P=29000; S=1450},
    @{Row=14; D=44427; M=200; N=20000; O=21000; P=20500; S=1025},
    @{Row=15; D=44420; M=160; N=20000; O=21000; P=20500; S=1025},
    @{Row=16; D=44874; M=240; N=29000; O=30000; P=29500; S=1475},
    @{Row=17; D=44473; M=40; N=19500; O=20000; P=19750; S=988},
    @{Row=18; D=44782; M=200; N=23500; O=24000; P=23750; S=1188},
    @{Row=19; D=44435; M=260; N=20000; O=22000; P=21115; S=1056},
    @{Row=20; D=44326; M=160; N=19500; O=20000; P=19750; S=988},
    @{Row=21; D=44441; M=160; N=20000; O=21000; P=20500; S=1025},
    @{Row=22; D=44882; M=120; N=28000; O=30000; P=29000; S=1450},
    @{Row=23; D=44781; M=160; N=23000; O=24000; P=23500; S=1175},
    @{Row=24; D=44776; M=160; N=23000; O=24000; P=23500; S=1175},
    @{Row=25; D=44417; M=160; N=20000; O=21000; P=20500; S=1025},
    @{Row=26; D=44809; M=60; N=27000; O=28000; P=27500; S=1375},
    @{Row=27; D=44315; M=100; N=20000; O=21000; P=20500; S=1025},
    @{Row=28; D=44365; M=100; N=20000; O=21000; P=20500; S=1025},
    @{Row=29; D=44778; M=100; N=23000; O=24000; P=23500; S=1175},
    @{Row=30; D=44448; M=100; N=20000; O=21000; P=20500; S=1025},
    @{Row=31; D=44335; M=200; N=19000; O=20000; P=19500; S=975},
    @{Row=32; D=44466; M=100; N=20000; O=21000; P=20500; S=1025},
    @{Row=33; D=44410; M=200; N=20000; O=21000; P=20500; S=1025},
    @{Row=34; D=44301; M=100; N=18000; O=19000; P=18500; S=925},
    @{Row=35; D=44431; M=160; N=21000; O=22000; P=21500; S=1075},
    @{Row=36; D=44810; M=100; N=27000; O=28000; P=27500; S=1375},
    @{Row=37; D=44407; M=160; N=20000; O=21000; P=20500; S=1025},
    @{Row=38; D=44336; M=100; N=19500; O=20000; P=19750; S=988},
    @{Row=39; D=44880; M=100; N=28000; O=30000; P=29000; S=1450},
    @{Row=40; D=44784; M=160; N=27000; O=28000; P=27500; S=1375},
    @{Row=41; D=44343; M=100; N=19500; O=20000; P=19750; S=988}
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 4).Value  = $u.D   # D - Fecha
    $ws.Cells.Item($r, 13).Value = $u.M   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $u.N   # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $u.O   # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $u.P   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $u.S   # S - Precio $/Kg
}
